$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (e.g. "146.30", "67.339.42").
# Force the cells to be treated as Text before writing so Excel does not silently
# coerce numeric-looking strings into real numbers (which would drop trailing
# zeros / alter the stored representation). Revert the style afterwards so the
# cells keep their original (default) style, matching the unedited cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.339.42"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.945.54"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "470.65"
$ws.Range("E5").Value = "  +7.89%  "
$ws.Range("D6").Value = "146.30"
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.732"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +9.18%  "
$ws.Range("D11").Value = "0.0000351"
$ws.Range("E11").Value = "  +9.54%  "
$ws.Range("D12").Value = "43.48"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "4.576.01"
$ws.Range("E13").Value = "  +4.18%  "
$ws.Range("D14").Value = "10.40"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "15.13"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "3.920.44"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "19.88"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "67.554.47"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "434.40"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("E22").Value = "  +4.87%  "
$ws.Range("D23").Value = "14.42"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("D24").Value = "87.56"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "3.61"
$ws.Range("E25").Value = "  +6.24%  "
$ws.Range("D26").Value = "38.71"
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "719.24"
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "13.48"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").Value = "42.30"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("D34").Value = "0.0₃0861"
$ws.Range("E34").Value = "  +25.09%  "
$ws.Range("D35").Value = "57.98"
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "5.35"
$ws.Range("E38").Value = "  -4.54%  "
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "3.51"
$ws.Range("E42").Value = "  +5.80%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.336"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  +6.84%  "
$ws.Range("D46").Value = "2.19"
$ws.Range("E46").Value = "  +5.05%  "
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  -7.40%  "
$ws.Range("D48").Value = "147.44"
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").Value = "3.19"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "2.88"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "25.73"
$ws.Range("E51").Value = "  +3.12%  "

# Restore the original (default) style for column D cells so no stray
# number-format style lingers on the cells in the saved file.
$ws.Range("D2:D51").Style = "Normal"
